$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.242.65"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "'1.901.41"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D5").Value = "'326.00"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.4639"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "'0.3915"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").Value = "'0.07883"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'0.9883"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'21.82"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("D12").Value = "'1.894.49"
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("D13").Value = "'7.074"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "'5.736"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "'0.06991"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "'88.24"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'0.000009972"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'17.11"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'29.244.40"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'2.153.33"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").Value = "'2.100"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").Value = "'155.96"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "'5.992"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Value = "'118.49"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "'1.882"
$ws.Range("E30").Value = "  -6.16%  "
$ws.Range("D31").Value = "'0.09353"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "'0.9006"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "'1.324"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").Value = "'3.212"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "'1.189"
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("D37").Value = "'0.05786"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").Value = "'0.02089"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "'1.001"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'7.711"
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("D41").Value = "'0.5715"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").Value = "'0.1787"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "'9.690"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").Value = "'11.96"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'0.5354"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "'2.176"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").Value = "'0.07026"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "'1.850"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").Value = "'2.567"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "'113.12"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'1.056"
$ws.Range("E51").Value = "  -1.40%  "
